$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F4").Value = "recall and mmr was signifactly imporved"
$ws.Range("C4").Value = "used the dictaber after the trankit was not efficent"
